$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.468.68'
$ws.Range("E2").Value = '  -1.70%  '
$ws.Range("D3").Value = '2.428.28'
$ws.Range("E3").Value = '  -2.25%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'511.17"
$ws.Range("E5").Value = '  -2.80%  '
$ws.Range("D6").Value = "'128.95"
$ws.Range("E6").Value = '  -3.50%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = "'0.547"
$ws.Range("E8").Value = '  -2.33%  '
$ws.Range("D9").Value = '2.438.50'
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("D11").Value = "'0.0946"
$ws.Range("E11").Value = '  -5.38%  '
$ws.Range("D12").Value = "'5.14"
$ws.Range("E12").Value = '  -5.32%  '
$ws.Range("E13").Value = '  -4.05%  '
$ws.Range("D14").Value = '2.858.97'
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("D15").Value = '57.384.65'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").Value = "'21.68"
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("E17").Value = '  -3.70%  '
$ws.Range("D18").Value = '2.434.18'
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("E19").Value = '  -4.92%  '
$ws.Range("D20").Value = "'314.59"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("E21").Value = '  -2.85%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").Value = '  -3.43%  '
$ws.Range("D24").Value = "'63.18"
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("E25").Value = '  -2.36%  '
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").Value = "'7.21"
$ws.Range("E28").Value = '  -3.56%  '
$ws.Range("D29").Value = "'169.88"
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").Value = "'6.20"
$ws.Range("E30").Value = '  -3.54%  '
$ws.Range("D31").Value = '0.0₃0716'
$ws.Range("E31").Value = '  -4.83%  '
$ws.Range("E32").Value = '  -2.89%  '
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("E36").Value = '  -3.46%  '
$ws.Range("E37").Value = '  -5.30%  '
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").Value = "'36.21"
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  -3.38%  '
$ws.Range("E41").Value = '  -3.84%  '
$ws.Range("D42").Value = "'271.21"
$ws.Range("E42").Value = '  -2.55%  '
$ws.Range("E43").Value = '  -4.83%  '
$ws.Range("E44").Value = '  -2.12%  '
$ws.Range("D45").Value = "'0.578"
$ws.Range("E45").Value = '  -2.85%  '
$ws.Range("D46").Value = "'0.0904"
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("D47").Value = "'119.86"
$ws.Range("E47").Value = '  -6.08%  '
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").Value = "'17.03"
$ws.Range("E49").Value = '  -4.73%  '
$ws.Range("E50").Value = '  -3.59%  '
$ws.Range("D51").Value = "'16.46"
$ws.Range("E51").Value = '  -4.58%  '
